$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 3 new rows (34-36) following the existing data pattern.
$data = @(
    @(10005, 110033, 10005),
    @(10005, 110034, 10005),
    @(10005, 110035, 10005)
)

$startRow = 34
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Update selection to mimic Excel selecting the entire row after the last data row
# (as if the row header for row 37 was clicked).
$ws.Range("A37:XFD1048576").Select()
